# Adds the "Packages:" section (handling of package creation) after the
# paragraph that ends with "...leicht erweitern kann. " and before the
# final (pre-existing) empty paragraph of the document body.
#
# Target structure (matches the XML diff):
#   <p>...leicht erweitern kann. </p>
#   <p></p>                                            <- new, empty
#   <p>Packages: </p>                                  <- new
#   <p>Ein Package besteht aus 5 Karten, ...</p>        <- new
#   <p></p>                                            <- pre-existing empty paragraph (34706413)

$d = $word.ActiveDocument

# Locate the pre-existing trailing empty paragraph (last paragraph of the body).
$trailing = $d.Paragraphs($d.Paragraphs.Count)

$packagesText = "Packages: "
$packageBodyText = "Ein Package besteht aus 5 Karten, diese können ident sein, bekommen aber unterschiedliche IDs. Damit man die Karten genau referenzieren kann, sonst wäre die Auswahl in der DB der entsprechenden Karte schwieriger. "

# Insert all three new paragraphs in one shot, just before the trailing
# empty paragraph. A leading marker character ("X") is used so the first
# (otherwise totally empty) paragraph is created with real run content —
# inserting a bare paragraph mark on its own leaves a stray placeholder
# run behind. We strip the marker off again right after.
$trailing.Range.InsertBefore("X`r" + $packagesText + "`r" + $packageBodyText + "`r")

# Remove the "X" marker character from the first of the newly inserted
# paragraphs, leaving it genuinely empty (no residual run).
$markerPara = $d.Paragraphs($d.Paragraphs.Count - 3)
$marker = $d.Range($markerPara.Range.Start, $markerPara.Range.Start + 1)
$marker.Delete()
